# Generate Report for Handback
# - Flip "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it appears
#   (Overview zh-cn/de-de status columns + per-language Status column).
# - Record the zh-cn handback timestamp (was the zero-date placeholder) and the
#   de-de handback timestamp (brand-new row).
# - Populate the "Latest Target File" (now a hyperlink, like column A) and
#   "Latest Handback File" columns for both language sheets.
# - Widen the columns that now hold longer content.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$targetFileName = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5b1ae90ef63d6d90d269548398fcfafecbd602e2/e2e/1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"

# ---- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ----
$newStatus = "Handed back: in sync with en-US"

$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# ---- zh-cn: Latest Target File / Latest Handback File / Latest Handback DateTime ----
$ws2.Hyperlinks.Add($ws2.Range("I2"), $targetUrl, "", "", $targetFileName)
$ws2.Hyperlinks.Add($ws2.Range("I3"), $targetUrl, "", "", $targetFileName)

$ws2.Range("J2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.zh-cn.xlf"
$ws2.Range("J3").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.zh-cn.xlf"

$ws2.Range("K2").Value = "2016-08-15 20:58:52"
$ws2.Range("K3").Value = "2016-08-15 20:58:52"

# ---- de-de: Latest Target File / Latest Handback File / Latest Handback DateTime ----
$ws3.Hyperlinks.Add($ws3.Range("I2"), $targetUrl, "", "", $targetFileName)
$ws3.Hyperlinks.Add($ws3.Range("I3"), $targetUrl, "", "", $targetFileName)

$ws3.Range("J2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.de-de.xlf"
$ws3.Range("J3").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.de-de.xlf"

$ws3.Range("K2").Value = "2016-08-15 20:59:00"
$ws3.Range("K3").Value = "2016-08-15 20:59:00"

# ---- Column widths ----
# Overview: Status columns (E, F) grew to fit the longer text.
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666664
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666664

# zh-cn / de-de: Status column (C) grew; Latest Target File / Latest Handback
# File columns (I, J) widened to fit filenames / URLs.
$ws2.Columns.Item(3).ColumnWidth = 29.166666666666664
$ws2.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws2.Columns.Item(10).ColumnWidth = 39.166666666666664

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666664
$ws3.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws3.Columns.Item(10).ColumnWidth = 39.166666666666664
